$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 493
$ws.Range("F6").Value = 62
$ws.Range("F7").Value = 1220
$ws.Range("F9").Value = 243
$ws.Range("F11").Value = 8205
$ws.Range("F13").Value = 9856
$ws.Range("F14").Value = 83
$ws.Range("F16").Value = 11
$ws.Range("F27").Value = 1691
$ws.Range("F33").Value = 547
$ws.Range("F38").Value = 416
$ws.Range("F40").Value = 9
$ws.Range("F41").Value = 121
$ws.Range("F43").Value = 309
$ws.Range("F45").Value = 257
$ws.Range("F46").Value = 104
$ws.Range("F48").Value = 20
$ws.Range("F49").Value = 20

# --- Sheet: 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 53
$ws.Range("F19").Value = 29
$ws.Range("F20").Value = 362

# --- Sheet: 本地生活 (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2763
$ws.Range("F5").Value = 191

# --- Sheet: 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 493
$ws.Range("F6").Value = 191
$ws.Range("F9").Value = 62
$ws.Range("F10").Value = 1220
$ws.Range("F14").Value = 243
$ws.Range("F16").Value = 8205
$ws.Range("F18").Value = 9856
$ws.Range("F19").Value = 11
$ws.Range("F24").Value = 1691
$ws.Range("F30").Value = 547
$ws.Range("F37").Value = 416
$ws.Range("F38").Value = 53
$ws.Range("F39").Value = 325
$ws.Range("F40").Value = 121
$ws.Range("F42").Value = 310
$ws.Range("F45").Value = 29
$ws.Range("F46").Value = 362
$ws.Range("F48").Value = 20
$ws.Range("F49").Value = 20
